# The commit swaps the presentation's theme from the "Integral" theme
# (ppt/theme/theme1.xml, bound to the slide master / whole deck) to the
# built-in PowerPoint "Office Theme" (whose colours previously lived,
# unused by any slide, in ppt/theme/theme2.xml bound only to the notes
# master). Font scheme / format scheme are identical between the two
# themes, so the only functional difference is the 12-colour theme
# colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's theme colour scheme indices (ThemeColorScheme.Item(n)),
# in order, are: dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4,
# accent5, accent6, hlink, folHlink. We push the "Office" theme's RGB
# values (converted to VBA's R + G*256 + B*65536 encoding) into each
# slot of the deck's theme colour scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$themeColors = $theme.ThemeColorScheme

# Best-effort: rename the theme / colour-scheme to match the built-in
# "Office Theme" naming (no-op on hosts where these are read-only).
try { $theme.Name = "Office Theme" } catch { }
try { $theme.ThemeElements.ColorScheme.Name = "Office" } catch { }
try { $theme.ThemeElements.FontScheme.Name = "Office" } catch { }

# Office Theme colours, in ThemeColorScheme index order 1..12:
#   1 dk1      = 000000
#   2 lt1      = FFFFFF
#   3 dk2      = 44546A
#   4 lt2      = E7E6E6
#   5 accent1  = 5B9BD5
#   6 accent2  = ED7D31
#   7 accent3  = A5A5A5
#   8 accent4  = FFC000
#   9 accent5  = 4472C4
#  10 accent6  = 70AD47
#  11 hlink    = 0563C1
#  12 folHlink = 954F72
$officeThemeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRgb[$i - 1]
}
